# The commit adds data to the (previously empty) worksheet: cell A1 gets
# the numeric value 34563 (a plain number, no special formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 34563
